$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9653424024581909
$ws.Range("B1").Value = 1.047032952308655
$ws.Range("C1").Value = 0.8838189244270325
$ws.Range("D1").Value = 0.9217614531517029
$ws.Range("E1").Value = 1.074640989303589
